$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting existing rows 58-72 down to 59-73
$ws.Rows(58).Insert()

# Populate the new row 58 with the new price record
$ws.Cells.Item(58, 1).Value = 7
$ws.Cells.Item(58, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(58, 3).Value = "Ñuble"
$ws.Cells.Item(58, 4).Value = 44637
$ws.Cells.Item(58, 5).Value = 16
$ws.Cells.Item(58, 6).Value = 100112021
$ws.Cells.Item(58, 7).Value = "Ají"
$ws.Cells.Item(58, 8).Value = "Americana (o)"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 80
$ws.Cells.Item(58, 11).Value = 8500
$ws.Cells.Item(58, 12).Value = 9000
$ws.Cells.Item(58, 13).Value = 8750
$ws.Cells.Item(58, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(58, 15).Value = "Región del Maule"
$ws.Cells.Item(58, 16).Value = 583
$ws.Cells.Item(58, 17).Value = 15
$ws.Cells.Item(58, 18).Value = "Hortaliza"
